$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update the mandataire's info and gross/net amount
$ws.Range("A2").Value = "JEMAA HORMI"
$ws.Range("B2").Value = "B219321"
$ws.Range("C2").Value = "'225400000805987601012173"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "KHOURIBGA"
$ws.Range("E2").Value = "CA"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "001/RRR"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 12000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 12000

# Row 3: previously blank filler row, now a second entry for the same person
$ws.Range("A3").Value = "JEMAA HORMI"
$ws.Range("B3").Value = "B219321"
$ws.Range("C3").Value = "'225400000805987601012173"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "KHOURIBGA"
$ws.Range("E3").Value = "CA"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "001/RRR"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2000

# Row 4: new blank filler row with the updated totals
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("I4").Value = 14000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 14000
